$wb = $excel.ActiveWorkbook

# New chloride-lab sampling rows, one appended to the bottom of each of the
# first nine data sheets (WIC, YS, SW, YI, YN, 6MC, DC, PBMS, PBSF).
# Column A = filter type ("Whatman", shared string already used throughout
# the workbook), B/C = datetime collected / datetime run (same timestamp -
# filter run immediately after collection). D (filtered date) and E
# (chloride result) are left blank - not measured/analyzed yet.

$rows = @(
    @{ Sheet = "WIC";  Row = 31; DateTime = 44264.453472222223 },
    @{ Sheet = "YS";   Row = 44; DateTime = 44264.469444444447 },
    @{ Sheet = "SW";   Row = 43; DateTime = 44264.490277777775 },
    @{ Sheet = "YI";   Row = 44; DateTime = 44264.4375 },
    @{ Sheet = "YN";   Row = 42; DateTime = 44264.54583333333 },
    @{ Sheet = "6MC";  Row = 44; DateTime = 44264.561111111114 },
    @{ Sheet = "DC";   Row = 44; DateTime = 44264.572222222225 },
    @{ Sheet = "PBMS"; Row = 44; DateTime = 44264.587500000001 },
    @{ Sheet = "PBSF"; Row = 43; DateTime = 44264.6 }
)

foreach ($entry in $rows) {
    $ws = $wb.Worksheets.Item($entry.Sheet)

    $rowNum = $entry.Row
    $prevRow = $rowNum - 1
    $dt = $entry.DateTime

    # PBSF's column A carries a one-off cell style (the "Good" cell style,
    # stripped back to plain formatting) on every existing filter-type
    # entry; carry that format down to the new row the way AutoFill/"extend
    # formatting" would, before writing the value.
    if ($entry.Sheet -eq "PBSF") {
        $ws.Range(("A" + $prevRow)).Copy()
        $ws.Range(("A" + $rowNum)).PasteSpecial(-4122)
    }

    $ws.Cells.Item($rowNum, 1).Value = "Whatman"
    $ws.Cells.Item($rowNum, 2).Value = $dt
    $ws.Cells.Item($rowNum, 3).Value = $dt

    $ws.Range(("A" + $prevRow + ":A" + $rowNum)).Select() | Out-Null
}

# Leave PBSF (the 9th / last-edited sheet) as the active tab and leave the
# selection on the newly-entered run-date cell, matching the original
# authoring session.
$wsLast = $wb.Worksheets.Item("PBSF")
$wsLast.Activate()
$wsLast.Range("C43").Select() | Out-Null
